$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.5902912020683289
$ws.Range("B1").Value = 1.030855774879456
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 1.837823152542114
$ws.Range("E1").Value = 1.561772108078003
